$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress percentages (column C)
$ws.Range("C16").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("C23").Value = 0.8
$ws.Range("C24").Value = 1
$ws.Range("C25").Value = 1

# Edit the comment text in D23: drop "après coup" from the sentence.
# (Re-assigning the text also moves this shared string to the end of the
# shared-strings table, same as the authored diff.)
$ws.Range("D23").Value = "utilisation de AmChart ou si pas de temps Voilà"

# Move the visible selection to match the saved view state.
$ws.Range("D24").Select()
